$d = $word.ActiveDocument

# 1. Title heading and the bold repeated title near the end (both occurrences
#    are replaced in one call since Replace:=2 is wdReplaceAll)
$d.Content.Find.Execute("Play Desperate Dawgs 2 Gigablox Free - Review 2021", $true, $false, $false, $false, $false, $true, 1, $false, "Play Desperate Dawgs 2 Gigablox Free: Game Review", 2) | Out-Null

# 2. "What we like" bullet points
$d.Content.Find.Execute("Exciting Gigablox mechanic for larger symbols", $true, $false, $false, $false, $false, $true, 1, $false, "Exciting gameplay mechanics with Gigablox symbols and expanding reels", 2) | Out-Null
$d.Content.Find.Execute("Three unique bonus features to keep the game interesting", $true, $false, $false, $false, $false, $true, 1, $false, "Three thrilling bonus features for added excitement and potential wins", 2) | Out-Null
$d.Content.Find.Execute("Golden Bet feature for increased chances of activating bonus", $true, $false, $false, $false, $false, $true, 1, $false, "Golden Bet feature to increase chances of activating the bonus round", 2) | Out-Null
$d.Content.Find.Execute("High potential for a sizable win", $true, $false, $false, $false, $false, $true, 1, $false, "High maximum win potential of up to 9,872 times the bet per spin", 2) | Out-Null

# 3. "What we don't like" bullet points
$d.Content.Find.Execute("RTP rate is slightly below industry average", $true, $false, $false, $false, $false, $true, 1, $false, "Slightly below-average RTP rate of 95.50%", 2) | Out-Null
$d.Content.Find.Execute("Limited variety in symbol design", $true, $false, $false, $false, $false, $true, 1, $false, "Limited variety of symbols on the reels", 2) | Out-Null

# 4. Final italic summary paragraph
$d.Content.Find.Execute("Read our review of Desperate Dawgs 2 Gigablox and play this exciting game for free today. Features the Gigablox mechanic for larger symbols and three unique bonus rounds.", $true, $false, $false, $false, $false, $true, 1, $false, "Discover the gameplay mechanics, features, and wins in Desperate Dawgs 2 Gigablox. Play for free now!", 2) | Out-Null
